# Actualización desde MV -datos-
# Appends 5 new daily rows (02-10-2021 .. 06-10-2021) after the last
# existing data row (275, "01-10-2021") on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("02-10-2021", 2087, 2604, 15060, 2629, 4552, 8745),
    @("03-10-2021", 2087, 2604, 15060, 2629, 4552, 8745),
    @("04-10-2021", 2108, 2630, 15210, 2655, 4597, 8832),
    @("05-10-2021", 2103, 2623, 15172, 2648, 4586, 8810),
    @("06-10-2021", 2090, 2608, 15084, 2633, 4559, 8759)
)

$startRow = 276
$endRow = $startRow + $newRows.Length - 1

# Pre-format column A for the new rows as Text so the dd-mm-yyyy-looking
# strings are stored verbatim (as shared strings) instead of being
# auto-converted into date serial numbers.
$ws.Range("A$startRow" + ":A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("E$r").Value = $data[4]
    $ws.Range("F$r").Value = $data[5]
    $ws.Range("G$r").Value = $data[6]
}

# Drop the temporary Text number format again so the new cells end up
# without any explicit style, matching the rest of the data rows.
$ws.Range("A$startRow" + ":A$endRow").ClearFormats()

Write-Host "Added rows $startRow to $endRow"
